# Update attendance status ("estado") from "sin dictar" to "asistio"
# for the rows belonging to ciclo 2 / semana 3 (3ra semana 2 ciclo).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(396,397,402,406,407,409,415,420,421,422,423,424,425,426,427,428,429,430,431,432,433,436,437,438,439,442,443,444,445,452)

foreach ($r in $rows) {
    $ws.Range("D$r").Value = "asistio"
}

# Reflect the final selected cell as recorded in the saved workbook.
$ws.Range("D578").Select()
